$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Students")

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "coco"

$ws.Range("B5").Select()
